# Daily attendance processing - swap the order of the last two names
# in the comma-separated "Recorded By" list (column G), leaving the
# literal "backup@backdoor.com, System" pairing untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-LastTwo($text) {
    if ($null -eq $text) {
        return $text
    }
    # Case-sensitive check: "backup@backdoor.com, System" is left untouched.
    # (PowerShell's -eq/-ceq are case-insensitive in this host, so use
    # the .NET string .Equals() method, which is case-sensitive by default.)
    if ($text.Equals("backup@backdoor.com, System")) {
        return $text
    }
    $parts = $text -split ", "
    if ($parts.Count -lt 2) {
        return $text
    }
    $n = $parts.Count
    $tmp = $parts[$n - 1]
    $parts[$n - 1] = $parts[$n - 2]
    $parts[$n - 2] = $tmp
    return ($parts -join ", ")
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    $updated = Swap-LastTwo $current
    if (-not $updated.Equals($current)) {
        $cell.Value2 = $updated
    }
}
